$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Rank" header in G1
$ws.Range("G1").Value = "Rank"

# Fill G2:G16 with sequential rank numbers 1..15
for ($i = 2; $i -le 16; $i++) {
    $ws.Cells.Item($i, 7).Value = $i - 1
}

# Update the view: scroll so row 9 is the top-left visible row, and select G16
$window = $excel.ActiveWindow
$window.ScrollRow = 9
$ws.Range("G16").Select()
